$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1829.9
$ws.Range("I40").Value = 1599.6
$ws.Range("K40").Value = 1599.6
$ws.Range("M40").Value = -1424.6
# Row 62
$ws.Range("H62").Value = 2366.6667
$ws.Range("I62").Value = 1233.3334
$ws.Range("J62").Value = 2933.3333
$ws.Range("K62").Value = 1233.3334
$ws.Range("L62").Value = 2933.3333
$ws.Range("M62").Value = -609.3334
$ws.Range("N62").Value = -4181.3333
# Row 65
$ws.Range("H65").Value = 2366.6667
$ws.Range("I65").Value = 1233.3334
$ws.Range("J65").Value = 2933.3333
$ws.Range("K65").Value = 6166.666999999999
$ws.Range("L65").Value = 14666.6665
$ws.Range("M65").Value = -3046.666999999999
$ws.Range("N65").Value = -20906.6665
# Row 98
$ws.Range("H98").Value = 529.6316
$ws.Range("I98").Value = 529.6316
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 529.6316
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 968.3684
$ws.Range("N98").ClearContents()
# Row 122
$ws.Range("H122").Value = 529.6316
$ws.Range("I122").Value = 529.6316
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1588.8948
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 861.1052
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1988.0526
$ws.Range("I132").Value = 1653.5491
$ws.Range("K132").Value = 4960.6473
$ws.Range("M132").Value = -2430.6473
# Row 138
$ws.Range("H138").Value = 2756.4783
$ws.Range("I138").Value = 1434.6666
$ws.Range("J138").Value = 3729.132
$ws.Range("K138").Value = 4303.9998
$ws.Range("L138").Value = 11187.396
$ws.Range("M138").Value = 836.0002000000004
$ws.Range("N138").Value = -21467.396

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 486866.66
$ws.Range("I32").Value = 524094.66
$ws.Range("J32").Value = 27721.334
$ws.Range("K32").Value = 524094.66
$ws.Range("L32").Value = 27721.334
$ws.Range("M32").Value = -523807.66
$ws.Range("N32").Value = -28295.334
# Row 102
$ws.Range("H102").Value = 4380
$ws.Range("I102").Value = 4631.4287
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 4631.4287
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -3009.4287
$ws.Range("N102").Value = -6744
# Row 132
$ws.Range("H132").Value = 3288.9583
$ws.Range("I132").Value = 2279.9375
$ws.Range("K132").Value = 6839.8125
$ws.Range("M132").Value = -4309.8125

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2813.1516
$ws.Range("I134").Value = 2277.0588
$ws.Range("J134").Value = 3382.75
$ws.Range("K134").Value = 6831.176399999999
$ws.Range("L134").Value = 10148.25
$ws.Range("M134").Value = -4296.176399999999
$ws.Range("N134").Value = -15218.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8666.25
$ws.Range("I31").Value = 1645.0714
$ws.Range("J31").Value = 13134.272
$ws.Range("K31").Value = 1645.0714
$ws.Range("L31").Value = 13134.272
$ws.Range("M31").Value = -1350.0714
$ws.Range("N31").Value = -13724.272
# Row 34
$ws.Range("H34").Value = 8666.25
$ws.Range("I34").Value = 1645.0714
$ws.Range("J34").Value = 13134.272
$ws.Range("K34").Value = 1645.0714
$ws.Range("L34").Value = 13134.272
$ws.Range("M34").Value = -1443.0714
$ws.Range("N34").Value = -13538.272
# Row 134
$ws.Range("H134").Value = 2765.7385
$ws.Range("I134").Value = 2654.3635
$ws.Range("K134").Value = 7963.0905
$ws.Range("M134").Value = -5428.0905
# Row 141
$ws.Range("H141").Value = 312500
$ws.Range("J141").Value = 350000
$ws.Range("L141").Value = 350000
$ws.Range("N141").Value = -360360

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 24001088
$ws.Range("J4").Value = 33335000
$ws.Range("L4").Value = 100005000
$ws.Range("N4").Value = -100005224
# Row 51
$ws.Range("H51").Value = 1204.762
$ws.Range("I51").Value = 700
$ws.Range("J51").Value = 1230
$ws.Range("K51").Value = 2100
$ws.Range("L51").Value = 3690
$ws.Range("M51").Value = -1640
$ws.Range("N51").Value = -4610
# Row 131
$ws.Range("H131").Value = 5286.5
$ws.Range("J131").Value = 6330.7393
$ws.Range("L131").Value = 18992.2179
$ws.Range("N131").Value = -29072.2179
# Row 132
$ws.Range("H132").Value = 2280.359
$ws.Range("J132").Value = 2354.8147
$ws.Range("L132").Value = 21193.3323
$ws.Range("N132").Value = -26253.3323
# Row 133
$ws.Range("H133").Value = 15962.143
$ws.Range("I133").Value = 6415
$ws.Range("J133").Value = 23122.5
$ws.Range("K133").Value = 19245
$ws.Range("L133").Value = 69367.5
$ws.Range("M133").Value = -14185
$ws.Range("N133").Value = -79487.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2799
$ws.Range("I126").Value = 2799
$ws.Range("K126").Value = 8397
$ws.Range("M126").Value = -5927

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 6035002
$ws.Range("J2").Value = 6035002
$ws.Range("L2").Value = 6035002
$ws.Range("N2").Value = -6035226
# Row 16
$ws.Range("H16").Value = 4609781
$ws.Range("I16").Value = 1419.28
$ws.Range("J16").Value = 23811288
$ws.Range("K16").Value = 1419.28
$ws.Range("L16").Value = 23811288
$ws.Range("M16").Value = -1249.28
$ws.Range("N16").Value = -23811628
# Row 62
$ws.Range("H62").Value = 90124.5
$ws.Range("J62").Value = 90124.5
$ws.Range("L62").Value = 90124.5
$ws.Range("N62").Value = -91372.5
# Row 65
$ws.Range("H65").Value = 90124.5
$ws.Range("J65").Value = 90124.5
$ws.Range("L65").Value = 270373.5
$ws.Range("N65").Value = -276613.5
# Row 69
$ws.Range("H69").Value = 34500
$ws.Range("J69").Value = 34500
$ws.Range("L69").Value = 34500
$ws.Range("N69").Value = -36122
# Row 72
$ws.Range("H72").Value = 34500
$ws.Range("J72").Value = 34500
$ws.Range("L72").Value = 103500
$ws.Range("N72").Value = -111612
# Row 128
$ws.Range("H128").Value = 46764.5
$ws.Range("J128").Value = 46764.5
$ws.Range("L128").Value = 46764.5
$ws.Range("N128").Value = -56724.5
# Row 132
$ws.Range("H132").Value = 2806.8857
$ws.Range("J132").Value = 2965.7856
$ws.Range("L132").Value = 8897.356800000001
$ws.Range("N132").Value = -13957.3568
# Row 139
$ws.Range("H139").Value = 963806.2
$ws.Range("I139").Value = 19300000
$ws.Range("J139").Value = 46996.5
$ws.Range("K139").Value = 19300000
$ws.Range("L139").Value = 46996.5
$ws.Range("M139").Value = -19294860
$ws.Range("N139").Value = -57276.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 113413780
$ws.Range("I2").Value = 5126000
$ws.Range("J2").Value = 200044000
$ws.Range("K2").Value = 5126000
$ws.Range("L2").Value = 200044000
$ws.Range("M2").Value = -5125888
$ws.Range("N2").Value = -200044224
